# "version final sin errores"
#
# Two logical changes to the "Metadata" sheet:
#   1. Bump the Version property value from 0.4.0 to 0.7.0 (row 3, col B).
#   2. Remove the Jurisdiction / Chile metadata row entirely (was row 11),
#      shifting every following row up by one.
#
# The "Concepts" sheet (sheet2) is untouched by this edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Update the Version value.
$ws.Range("B3").Value = "0.7.0"

# 2. Delete the entire "Jurisdiction" / "Chile" row (row 11), shifting
#    subsequent rows up so the sheet ends up as A1:B21.
$ws.Rows.Item(11).Delete()
